$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "sec8lec65 czemu switch musi mieć wartość default skoro przeskanowalismy wszystkie wartości? Czy java wie jaki jest pełen zakres danych do sprawdzenia?"

$ws.Range("C16").Select()
